$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.535.44'
$ws.Cells.Item(2, 5).Value = '  +1.17%  '

$ws.Cells.Item(3, 4).Value = '1.773.27'
$ws.Cells.Item(3, 5).Value = '  -0.65%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.012'
$ws.Cells.Item(4, 5).Value = '  +0.77%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '337.43'
$ws.Cells.Item(5, 5).Value = '  +0.69%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.007'
$ws.Cells.Item(6, 5).Value = '  +0.58%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.3872'
$ws.Cells.Item(7, 5).Value = '  +2.48%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3426'
$ws.Cells.Item(8, 5).Value = '  -0.26%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '47.27'
$ws.Cells.Item(9, 5).Value = '  -2.14%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '1.152'
$ws.Cells.Item(10, 5).Value = '  -3.69%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07487'
$ws.Cells.Item(11, 5).Value = '  +0.17%  '

$ws.Cells.Item(12, 2).Value = 'BinanceUSD'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.007'
$ws.Cells.Item(12, 5).Value = '  +0.56%  '

$ws.Cells.Item(13, 2).Value = 'Solana'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '22.49'
$ws.Cells.Item(13, 5).Value = '  +3.64%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.413'
$ws.Cells.Item(14, 5).Value = '  -0.73%  '

$ws.Cells.Item(15, 2).Value = 'WrappedEther'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(15, 4).Value = '1.780.31'
$ws.Cells.Item(15, 5).Value = '  -0.51%  '

$ws.Cells.Item(16, 2).Value = 'Chainlink'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '7.145'
$ws.Cells.Item(16, 5).Value = '  +0.83%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.00001080'
$ws.Cells.Item(17, 5).Value = '  -1.37%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.06719'
$ws.Cells.Item(18, 5).Value = '  +0.80%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '83.24'
$ws.Cells.Item(19, 5).Value = '  -0.78%  '

$ws.Cells.Item(20, 5).Value = '  +0.42%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '17.57'
$ws.Cells.Item(21, 5).Value = '  +1.66%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.476'
$ws.Cells.Item(22, 5).Value = '  -2.13%  '

$ws.Cells.Item(23, 4).Value = '27.569.24'
$ws.Cells.Item(23, 5).Value = '  +1.24%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '12.23'
$ws.Cells.Item(24, 5).Value = '  -1.28%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.386'
$ws.Cells.Item(25, 5).Value = '  -1.50%  '

$ws.Cells.Item(26, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.480'
$ws.Cells.Item(26, 5).Value = '  -2.33%  '

$ws.Cells.Item(27, 2).Value = 'EthereumClassic'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '20.91'
$ws.Cells.Item(27, 5).Value = '  -2.12%  '

$ws.Cells.Item(28, 2).Value = 'ImmutableX'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.436'
$ws.Cells.Item(28, 5).Value = '  -4.60%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '154.35'
$ws.Cells.Item(29, 5).Value = '  +0.52%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '135.77'
$ws.Cells.Item(30, 5).Value = '  +1.27%  '

$ws.Cells.Item(31, 4).Value = '1.981.12'
$ws.Cells.Item(31, 5).Value = '  -0.42%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '6.202'
$ws.Cells.Item(32, 5).Value = '  +1.94%  '

$ws.Cells.Item(33, 5).Value = '  -0.76%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.08928'
$ws.Cells.Item(34, 5).Value = '  +2.84%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '12.89'
$ws.Cells.Item(35, 5).Value = '  -2.46%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.02451'
$ws.Cells.Item(36, 5).Value = '  +4.82%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '5.433'
$ws.Cells.Item(37, 5).Value = '  -0.48%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.6831'
$ws.Cells.Item(38, 5).Value = '  -1.52%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.06414'
$ws.Cells.Item(39, 5).Value = '  +1.35%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.2211'
$ws.Cells.Item(40, 5).Value = '  +0.65%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.533'
$ws.Cells.Item(41, 5).Value = '  -7.57%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.250'
$ws.Cells.Item(42, 5).Value = '  +0.64%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '8.448'
$ws.Cells.Item(43, 5).Value = '  -3.75%  '

$ws.Cells.Item(44, 5).Value = '  -0.29%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.005'
$ws.Cells.Item(45, 5).Value = '  +0.49%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.6335'
$ws.Cells.Item(46, 5).Value = '  -2.46%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.858'
$ws.Cells.Item(47, 5).Value = '  +0.19%  '

$ws.Cells.Item(48, 2).Value = 'Quant'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '132.48'
$ws.Cells.Item(48, 5).Value = '  +2.65%  '

$ws.Cells.Item(49, 2).Value = 'NEARProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.121'
$ws.Cells.Item(49, 5).Value = '  -0.98%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.07428'
$ws.Cells.Item(50, 5).Value = '  +4.23%  '

$ws.Cells.Item(51, 2).Value = 'EOS'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.160'
$ws.Cells.Item(51, 5).Value = '  +3.80%  '
